$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="37.397.24"'
$ws.Range("E2").Formula = '="  +2.04%  "'
$ws.Range("D3").Formula = '="2.058.87"'
$ws.Range("E3").Formula = '="  +1.86%  "'
$ws.Range("E4").Formula = '="  +0.13%  "'
$ws.Range("D5").Formula = '="233.21"'
$ws.Range("E5").Formula = '="  -0.81%  "'
$ws.Range("E6").Formula = '="  +3.13%  "'
$ws.Range("B7").Formula = '="USDC"'
$ws.Range("C7").Formula = '="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"'
$ws.Range("D7").Formula = '="1.00"'
$ws.Range("E7").Formula = '="  -0.06%  "'
$ws.Range("B8").Formula = '="Solana"'
$ws.Range("C8").Formula = '="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"'
$ws.Range("D8").Formula = '="57.68"'
$ws.Range("E8").Formula = '="  +5.12%  "'
$ws.Range("D9").Formula = '="0.383"'
$ws.Range("E9").Formula = '="  +3.44%  "'
$ws.Range("D10").Formula = '="58.73"'
$ws.Range("E10").Formula = '="  +1.24%  "'
$ws.Range("D11").Formula = '="0.0761"'
$ws.Range("E11").Formula = '="  +1.62%  "'
$ws.Range("E12").Formula = '="  +0.02%  "'
$ws.Range("D13").Formula = '="2.367.16"'
$ws.Range("E13").Formula = '="  +2.07%  "'
$ws.Range("D14").Formula = '="14.36"'
$ws.Range("E14").Formula = '="  +1.47%  "'
$ws.Range("D15").Formula = '="20.85"'
$ws.Range("E15").Formula = '="  +3.46%  "'
$ws.Range("D16").Formula = '="0.772"'
$ws.Range("E16").Formula = '="  +1.21%  "'
$ws.Range("D17").Formula = '="5.17"'
$ws.Range("E17").Formula = '="  +1.23%  "'
$ws.Range("D18").Formula = '="2.063.17"'
$ws.Range("E18").Formula = '="  +2.04%  "'
$ws.Range("D19").Formula = '="37.593.94"'
$ws.Range("E19").Formula = '="  +4.18%  "'
$ws.Range("D20").Formula = '="6.15"'
$ws.Range("E20").Formula = '="  +14.82%  "'
$ws.Range("D21").Formula = '="69.15"'
$ws.Range("E21").Formula = '="  +1.99%  "'
$ws.Range("D22").Formula = '="0.0₃0811"'
$ws.Range("E22").Formula = '="  +1.70%  "'
$ws.Range("D23").Formula = '="225.17"'
$ws.Range("E23").Formula = '="  +2.20%  "'
$ws.Range("E24").Formula = '="  +0.03%  "'
$ws.Range("D25").Formula = '="2.42"'
$ws.Range("E25").Formula = '="  +0.94%  "'
$ws.Range("D26").Formula = '="2.39"'
$ws.Range("E26").Formula = '="  +0.21%  "'
$ws.Range("D27").Formula = '="165.68"'
$ws.Range("E27").Formula = '="  +1.27%  "'
$ws.Range("D28").Formula = '="1.48"'
$ws.Range("E28").Formula = '="  +6.70%  "'
$ws.Range("D29").Formula = '="8.91"'
$ws.Range("E29").Formula = '="  +3.18%  "'
$ws.Range("D30").Formula = '="0.129"'
$ws.Range("E30").Formula = '="  +0.38%  "'
$ws.Range("D31").Formula = '="19.07"'
$ws.Range("E31").Formula = '="  +0.84%  "'
$ws.Range("E32").Formula = '="  +0.51%  "'
$ws.Range("D33").Formula = '="4.47"'
$ws.Range("E33").Formula = '="  +2.52%  "'
$ws.Range("B34").Formula = '="Hedera"'
$ws.Range("C34").Formula = '="https://coinranking.com/coin/jad286TjB+hedera-hbar"'
$ws.Range("D34").Formula = '="0.0620"'
$ws.Range("E34").Formula = '="  +2.74%  "'
$ws.Range("B35").Formula = '="LidoDAOToken"'
$ws.Range("C35").Formula = '="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"'
$ws.Range("D35").Formula = '="2.57"'
$ws.Range("E35").Formula = '="  +5.18%  "'
$ws.Range("D36").Formula = '="4.57"'
$ws.Range("E36").Formula = '="  +7.46%  "'
$ws.Range("B37").Formula = '="BinanceUSD"'
$ws.Range("C37").Formula = '="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"'
$ws.Range("D37").Formula = '="1.00"'
$ws.Range("E37").Formula = '="  +0.30%  "'
$ws.Range("B38").Formula = '="THORChain"'
$ws.Range("C38").Formula = '="https://coinranking.com/coin/ybmU-kKU+thorchain-rune"'
$ws.Range("D38").Formula = '="5.97"'
$ws.Range("E38").Formula = '="  +4.60%  "'
$ws.Range("B39").Formula = '="WEMIXToken"'
$ws.Range("C39").Formula = '="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"'
$ws.Range("D39").Formula = '="1.76"'
$ws.Range("E39").Formula = '="  -0.83%  "'
$ws.Range("B40").Formula = '="RenderToken"'
$ws.Range("C40").Formula = '="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"'
$ws.Range("D40").Formula = '="3.30"'
$ws.Range("E40").Formula = '="  +0.37%  "'
$ws.Range("D41").Formula = '="4.67"'
$ws.Range("E41").Formula = '="  +10.76%  "'
$ws.Range("E42").Formula = '="  -0.79%  "'
$ws.Range("D43").Formula = '="0.0947"'
$ws.Range("E43").Formula = '="  +2.07%  "'
$ws.Range("D44").Formula = '="1.454.12"'
$ws.Range("E44").Formula = '="  -0.26%  "'
$ws.Range("D45").Formula = '="1.17"'
$ws.Range("E45").Formula = '="  +5.95%  "'
$ws.Range("D46").Formula = '="95.28"'
$ws.Range("E46").Formula = '="  +5.84%  "'
$ws.Range("D47").Formula = '="0.0211"'
$ws.Range("E47").Formula = '="  +3.59%  "'
$ws.Range("D48").Formula = '="15.63"'
$ws.Range("E48").Formula = '="  +2.03%  "'
$ws.Range("E49").Formula = '="  +2.08%  "'
$ws.Range("D50").Formula = '="7.17"'
$ws.Range("E50").Formula = '="  +4.34%  "'
$ws.Range("E51").Formula = '="  +2.33%  "'

$dataRange = $ws.Range("A1:E51")
$dataRange.Copy() | Out-Null
$dataRange.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
